# Applies the "latest game results" update described by the commit:
#   - 4 new games (results rows 114-126 on "results", rows 34-37 on "games")
#   - refreshed sheet view / active sheet state left by the author's last save
#
# NOTE on write order: the shared-strings table is append-on-first-use, and
# the target OOXML has the 4 new game ids (indices 112-115) created before
# the 4 new game names (indices 116-119, which themselves land in REVERSE
# row order - row 37's name first, then 36, 35, 34). We reproduce that exact
# ordering below so xl/sharedStrings.xml comes out byte-identical in layout.

$wb = $excel.ActiveWorkbook
$results = $wb.Worksheets.Item("results")
$games = $wb.Worksheets.Item("games")

# ---------------------------------------------------------------------
# Step 1: "results" sheet - append the per-player rows for the 4 new games.
# Writing this sheet first mints shared strings 112-115 (the game ids) in
# game order, matching the diff.
# ---------------------------------------------------------------------
$results.Cells.Item(114, 1).Value = "694079731cb455091b9a92cb"
$results.Cells.Item(114, 2).Value = "Thomas"
$results.Cells.Item(114, 3).Value = 10
$results.Cells.Item(114, 4).Value = "Empyrean"

$results.Cells.Item(115, 1).Value = "694079731cb455091b9a92cb"
$results.Cells.Item(115, 2).Value = "Eric"
$results.Cells.Item(115, 3).Value = 8
$results.Cells.Item(115, 4).Value = "Barony"

$results.Cells.Item(116, 1).Value = "694079731cb455091b9a92cb"
$results.Cells.Item(116, 2).Value = "Manu"
$results.Cells.Item(116, 3).Value = 6
$results.Cells.Item(116, 4).Value = "Jol"

$results.Cells.Item(117, 1).Value = "69505c9be44f36bfafb98b45"
$results.Cells.Item(117, 2).Value = "Eric"
$results.Cells.Item(117, 3).Value = 10
$results.Cells.Item(117, 4).Value = "Sol"

$results.Cells.Item(118, 1).Value = "69505c9be44f36bfafb98b45"
$results.Cells.Item(118, 2).Value = "Manu"
$results.Cells.Item(118, 3).Value = 7
$results.Cells.Item(118, 4).Value = "Muaat"

$results.Cells.Item(119, 1).Value = "69505c9be44f36bfafb98b45"
$results.Cells.Item(119, 2).Value = "Thomas"
$results.Cells.Item(119, 3).Value = 6
$results.Cells.Item(119, 4).Value = "Arborec"

$results.Cells.Item(120, 1).Value = "6962727b1e677f7a6f600b1e"
$results.Cells.Item(120, 2).Value = "Thomas"
$results.Cells.Item(120, 3).Value = 10
$results.Cells.Item(120, 4).Value = "Yssaril"

$results.Cells.Item(121, 1).Value = "6962727b1e677f7a6f600b1e"
$results.Cells.Item(121, 2).Value = "Manu"
$results.Cells.Item(121, 3).Value = 8
$results.Cells.Item(121, 4).Value = "Arborec"

$results.Cells.Item(122, 1).Value = "6962727b1e677f7a6f600b1e"
$results.Cells.Item(122, 2).Value = "Eric"
$results.Cells.Item(122, 3).Value = 7
$results.Cells.Item(122, 4).Value = "Nekro"

$results.Cells.Item(123, 1).Value = "696cdb2972af55390d4e54d5"
$results.Cells.Item(123, 2).Value = "Manu"
$results.Cells.Item(123, 3).Value = 10
$results.Cells.Item(123, 4).Value = "Hacan"

$results.Cells.Item(124, 1).Value = "696cdb2972af55390d4e54d5"
$results.Cells.Item(124, 2).Value = "Thomas"
$results.Cells.Item(124, 3).Value = 9
$results.Cells.Item(124, 4).Value = "Winnu"

$results.Cells.Item(125, 1).Value = "696cdb2972af55390d4e54d5"
$results.Cells.Item(125, 2).Value = "Frank"
$results.Cells.Item(125, 3).Value = 8
$results.Cells.Item(125, 4).Value = "Barony"

$results.Cells.Item(126, 1).Value = "696cdb2972af55390d4e54d5"
$results.Cells.Item(126, 2).Value = "Eric"
$results.Cells.Item(126, 3).Value = 4
$results.Cells.Item(126, 4).Value = "Sardakk"

# ---------------------------------------------------------------------
# Step 2: "games" sheet - game ids (column A) for the 4 new games, top to
# bottom. These ids already exist in the shared-strings table from step 1,
# so this reuses indices 112-115 without minting new entries.
# ---------------------------------------------------------------------
$games.Cells.Item(34, 1).Value = "694079731cb455091b9a92cb"
$games.Cells.Item(35, 1).Value = "69505c9be44f36bfafb98b45"
$games.Cells.Item(36, 1).Value = "6962727b1e677f7a6f600b1e"
$games.Cells.Item(37, 1).Value = "696cdb2972af55390d4e54d5"

# ---------------------------------------------------------------------
# Step 3: "games" sheet - game names (column B), entered bottom row first
# (row 37) up to the top (row 34), which is the order the author actually
# typed them in and mints shared strings 116-119 in that same order.
# ---------------------------------------------------------------------
$games.Cells.Item(37, 2).Value = "TI We verlängert, wenn ihr euch noch in die Augen schauen könnt"
$games.Cells.Item(36, 2).Value = "Jetzt mal ohne Fehler bitte"
$games.Cells.Item(35, 2).Value = "Tommy ist einfach Tommy "
$games.Cells.Item(34, 2).Value = "The return of the Timmy"

# ---------------------------------------------------------------------
# Step 4: "games" sheet - remaining numeric/date columns for the 4 new rows.
# ---------------------------------------------------------------------
$games.Cells.Item(34, 3).Value = 10
$games.Cells.Item(34, 4).Value = 46006
$games.Cells.Item(34, 5).Value = 46018
$games.Cells.Item(34, 6).Value = 5

$games.Cells.Item(35, 3).Value = 10
$games.Cells.Item(35, 4).Value = 46018
$games.Cells.Item(35, 5).Value = 46032
$games.Cells.Item(35, 6).Value = 5

$games.Cells.Item(36, 3).Value = 10
$games.Cells.Item(36, 4).Value = 46032
$games.Cells.Item(36, 5).Value = 46037
$games.Cells.Item(36, 6).Value = 6

$games.Cells.Item(37, 3).Value = 10
$games.Cells.Item(37, 4).Value = 46040
$games.Cells.Item(37, 5).Value = 46055
$games.Cells.Item(37, 6).Value = 5

# ---------------------------------------------------------------------
# Step 5: leave the view the way the author's last save left it - "games"
# as the active/selected sheet, scrolled near the bottom with F35 selected;
# "results" scrolled down to the newly-added block with A114:XFD117 selected.
# ---------------------------------------------------------------------
[void]$results.Activate()
[void]$results.Range("A114:XFD117").Select()

[void]$games.Activate()
[void]$games.Range("F35").Select()
